$wb = $excel.ActiveWorkbook

# --- Add the new "inputTable" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "inputTable"

# --- Header row ---
$ws.Range("A1").Value = "Tower"
$ws.Range("B1").Value = "Panel"
$ws.Range("C1").Value = "Bracing"

# --- Tower 1 / Panels 1-7 -> singleCross ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "singleCross"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "singleCross"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "singleCross"

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "singleCross"

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "singleCross"

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = "singleCross"

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "singleCross"

# --- Separator row ---
$ws.Range("B9").Value = [char]0x2026

# --- Tower 2 / Panels 1-7 -> doubleCross ---
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "doubleCross"

$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "doubleCross"

$ws.Range("A12").Value = 2
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "doubleCross"

$ws.Range("A13").Value = 2
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = "doubleCross"

$ws.Range("A14").Value = 2
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = "doubleCross"

$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 6
$ws.Range("C15").Value = "doubleCross"

$ws.Range("A16").Value = 2
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = "doubleCross"

# --- Separator row ---
$ws.Range("B17").Value = [char]0x2026

# --- Make the new sheet the active / selected tab, with A1:C17 selected ---
$ws.Select() | Out-Null
$ws.Range("A1:C17").Select() | Out-Null
